# Applies the "poster updated.pptx" edits:
#  - reposition/resize the "Oval 38" circle (made smaller and nudged)
#  - reposition the "Oval 37" circle
#  - reposition the "Oval 19" circle (x-only nudge)
#  - reposition the "TextBox 7" (CKD intro) textbox (x-only nudge)
#  - reposition the "TextBox 14" (Conclusion ...) textbox (x-only nudge)
#  - collapse three runs of "and "/"Recommendations "/"for Further Study"
#    into a single run of text
#  - fix a typo: "Jiang" -> "Jian" in the acknowledgements paragraph

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---- Shape: Oval 38 (small circle, top-left area) ----
$oval38 = $s.Shapes.Item(1)
$oval38.Left   = 35.49913386
$oval38.Top    = 461.5469291
$oval38.Width  = 224.4819685
$oval38.Height = 234.4752756

# ---- Shape: Oval 37 (small circle, lower-left area) ----
$oval37 = $s.Shapes.Item(2)
$oval37.Left = 50.99960630
$oval37.Top  = 1737.395591

# ---- Shape: Oval 19 (small circle, right-hand column) ----
$oval19 = $s.Shapes.Item(4)
$oval19.Left = 2581.007165

# ---- Shape: TextBox 7 ("CKD" intro paragraph) ----
$textBox7 = $s.Shapes.Item(6)
$textBox7.Left = 87.72795276

# ---- Shape: TextBox 14 ("Conclusion and Recommendations ..." paragraph) ----
$textBox14 = $s.Shapes.Item(11)
$textBox14.Left = 2616.112913

# ---- Text edit: merge "and " / "Recommendations " / "for Further Study" ----
$tr14 = $textBox14.TextFrame.TextRange
$heading = $tr14.Find("and Recommendations for Further Study")
$heading.Text = "and Recommendations for Further Study"

# ---- Text edit: fix typo "Jiang" -> "Jian" in acknowledgements ----
$typo = $tr14.Find("Jiang")
$typo.Text = "Jian"

Write-Host "Applied poster layout + text fixes"
